# Apply updates to the three timetable sheets: Regular_Timetable,
# PreMid_Timetable, PostMid_Timetable.
#
# Changes per sheet:
#   B2  : "MINOR: Generative Ai [C101]" -> "MINOR: Generative Ai [C102]"
#   C8  : "DS302 (Lab) [L207]"          -> "DS302 (Lab) [L107]"
#   C9  : "DS302 (Lab) [L207]"          -> "DS302 (Lab) [L107]"
#   B10 : "MINOR: VLSI [C101]"          -> "MINOR: VLSI [C102]"

$wb = $excel.ActiveWorkbook

$sheetNames = @("Regular_Timetable", "PreMid_Timetable", "PostMid_Timetable")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
    $ws.Range("C8").Value = "DS302 (Lab) [L107]"
    $ws.Range("C9").Value = "DS302 (Lab) [L107]"
    $ws.Range("B10").Value = "MINOR: VLSI [C102]"
}
